# Generate Report for Handoff
#
# A new handoff run produced a fresh source-file GUID and xliff content
# hash, plus later timestamps for the same handoff. Refresh the
# localization-status report so every sheet (Overview, zh-cn, de-de)
# reflects the new run while leaving everything else (links' actual
# target URLs, table definitions, styles, ...) untouched.

$wb = $excel.ActiveWorkbook

$oldGuid = "be5e6cf7-9782-4ed4-adc6-fcace5b6fe73"
$newGuid = "236f1989-d7a8-4e4b-9869-e1217105f4ec"

$oldHash = "375d9e82c607793dae34eb7f69910b23d2f2af54"
$newHash = "cc18edd51dd83142b3f30b9c4daea148e9e3d718"

# Same external link target on every sheet - only the visible display
# text changes, the hyperlinked URL itself is untouched.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/oltest/blob/c19fe35f8ea9024743d4bd10c9d77b2577601077/e2e/$oldGuid.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-12 15:09:49"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-12 15:09:41"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-12 15:09:49"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
